$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2025-01-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-31 Friday", 2) | Out-Null

# Update each arithmetic-problem cell in the table, in row-major order
$tbl = $d.Tables.Item(1)
$values = @(
    "14+47=",
    "93-5=",
    "67-19=",
    "11-0=",
    "75-21=",
    "74-60=",
    "95-89=",
    "97-24=",
    "34-28=",
    "16+81=",
    "1+47=",
    "75-42=",
    "33-5=",
    "60-35=",
    "53+28=",
    "24-11=",
    "31+0=",
    "29-14=",
    "44-44=",
    "67+10=",
    "79+5=",
    "11+33=",
    "14+4=",
    "33+42=",
    "96-29=",
    "3+94=",
    "71+15=",
    "51+10=",
    "13+71=",
    "61+23=",
    "75+8=",
    "1+33=",
    "87-73=",
    "87-59=",
    "7-0=",
    "46-22=",
    "48+35=",
    "97-32=",
    "19+24=",
    "97-63=",
    "9+59=",
    "37+30=",
    "76-73=",
    "42-13=",
    "11+80=",
    "18+5=",
    "9+40=",
    "24+20=",
    "47+33=",
    "14+15=",
    "21+21=",
    "75-56=",
    "98-1=",
    "92-18=",
    "77-45=",
    "4+83=",
    "16+9=",
    "63+9=",
    "83-74=",
    "54+20=",
    "91-89=",
    "95-43=",
    "18-14=",
    "93-24=",
    "99-17=",
    "40-24=",
    "23+72=",
    "36+9=",
    "74+14=",
    "86-27=",
    "11+18=",
    "37+17=",
    "98-27=",
    "43+23=",
    "15+21=",
    "47-9=",
    "28+56=",
    "77-52=",
    "81+0=",
    "27+65=",
    "92-40=",
    "68-51=",
    "56+0=",
    "19+67=",
    "98-6=",
    "89-64=",
    "44+26=",
    "83-62=",
    "58-56=",
    "97-14=",
    "33+18=",
    "11+17=",
    "67-9=",
    "24+10=",
    "35+63=",
    "62+11=",
    "46+11=",
    "69+8=",
    "85-70=",
    "0+91="
)

$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Output "done: updated $idx cells"
